# Auto-generated edit script applying the cell-level price/volume refresh
# described in the commit diff for cryptos.xlsx. A handful of rows were
# also reordered (their Coin/Link/Price/Volume cells swapped with an
# adjacent row), which is reproduced below as straightforward per-cell
# value updates keyed by final cell reference.
#
# Every cell is forced to a Text number format before the value is
# written (and the format is cleared again afterwards) so that
# numeric-looking strings such as "0.999" or "185.44" are stored as
# text, matching the original inline-string cell type instead of being
# auto-converted to a number by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '74.744.17'
Set-TextValue 'E2' '  +8.78%  '
Set-TextValue 'D3' '2.592.64'
Set-TextValue 'E3' '  +6.57%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '185.44'
Set-TextValue 'E5' '  +15.32%  '
Set-TextValue 'D6' '580.65'
Set-TextValue 'E6' '  +3.84%  '
Set-TextValue 'D7' '0.999'
Set-TextValue 'B8' 'Dogecoin'
Set-TextValue 'C8' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D8' '0.207'
Set-TextValue 'E8' '  +25.61%  '
Set-TextValue 'B9' 'XRP'
Set-TextValue 'C9' 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue 'D9' '0.533'
Set-TextValue 'E9' '  +4.66%  '
Set-TextValue 'D10' '2.592.07'
Set-TextValue 'E10' '  +6.71%  '
Set-TextValue 'E11' '  -0.25%  '
Set-TextValue 'E12' '  +8.63%  '
Set-TextValue 'D13' '4.78'
Set-TextValue 'E13' '  +3.71%  '
Set-TextValue 'D14' '0.0000192'
Set-TextValue 'E14' '  +10.01%  '
Set-TextValue 'D15' '74.523.91'
Set-TextValue 'E15' '  +8.71%  '
Set-TextValue 'D17' '26.21'
Set-TextValue 'E17' '  +13.20%  '
Set-TextValue 'D18' '2.584.10'
Set-TextValue 'E18' '  +6.41%  '
Set-TextValue 'D19' '8.87'
Set-TextValue 'E19' '  +28.58%  '
Set-TextValue 'D20' '11.78'
Set-TextValue 'E20' '  +12.23%  '
Set-TextValue 'D21' '376.91'
Set-TextValue 'E21' '  +12.21%  '
Set-TextValue 'D22' '2.32'
Set-TextValue 'E22' '  +20.72%  '
Set-TextValue 'D23' '4.06'
Set-TextValue 'E23' '  +6.09%  '
Set-TextValue 'E24' '  +0.04%  '
Set-TextValue 'D25' '69.82'
Set-TextValue 'E25' '  +4.49%  '
Set-TextValue 'D26' '4.17'
Set-TextValue 'E26' '  +13.34%  '
Set-TextValue 'D27' '9.25'
Set-TextValue 'E27' '  +12.49%  '
Set-TextValue 'E28' '  +6.36%  '
Set-TextValue 'D29' '0.996'
Set-TextValue 'E29' '  -0.82%  '
Set-TextValue 'D30' '0.0₃0945'
Set-TextValue 'E30' '  +15.17%  '
Set-TextValue 'B31' 'Bittensor'
Set-TextValue 'C31' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D31' '511.47'
Set-TextValue 'E31' '  +19.37%  '
Set-TextValue 'B32' 'InternetComputer(DFINITY)'
Set-TextValue 'C32' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D32' '7.96'
Set-TextValue 'E32' '  +11.31%  '
Set-TextValue 'D33' '1.36'
Set-TextValue 'E33' '  +19.12%  '
Set-TextValue 'E34' '  +6.85%  '
Set-TextValue 'E35' '  -0.02%  '
Set-TextValue 'D36' '0.119'
Set-TextValue 'E36' '  +12.41%  '
Set-TextValue 'D37' '159.49'
Set-TextValue 'E37' '  -0.79%  '
Set-TextValue 'E38' '  +6.86%  '
Set-TextValue 'D39' '19.41'
Set-TextValue 'E39' '  +1.81%  '
Set-TextValue 'E40' '  -0.04%  '
Set-TextValue 'D41' '4.93'
Set-TextValue 'E41' '  +13.45%  '
Set-TextValue 'D42' '1.68'
Set-TextValue 'E42' '  +12.22%  '
Set-TextValue 'D43' '0.322'
Set-TextValue 'E43' '  +8.03%  '
Set-TextValue 'D44' '2.46'
Set-TextValue 'E44' '  +19.99%  '
Set-TextValue 'B45' 'ImmutableX'
Set-TextValue 'C45' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D45' '1.17'
Set-TextValue 'E45' '  +8.41%  '
Set-TextValue 'B46' 'OKB'
Set-TextValue 'C46' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D46' '38.89'
Set-TextValue 'E46' '  +4.18%  '
Set-TextValue 'D47' '151.98'
Set-TextValue 'E47' '  +15.55%  '
Set-TextValue 'D48' '0.0822'
Set-TextValue 'E48' '  +15.01%  '
Set-TextValue 'D49' '3.63'
Set-TextValue 'E49' '  +8.40%  '
Set-TextValue 'D50' '0.522'
Set-TextValue 'E50' '  +8.09%  '
Set-TextValue 'B51' 'Mantle'
Set-TextValue 'C51' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D51' '0.582'
Set-TextValue 'E51' '  +4.43%  '
